$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Fill in the remaining turnip price cells for week 1 (row 4)
$ws.Range("I4").Value = 76
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 78
$ws.Range("L4").Value = 54
$ws.Range("M4").Value = 65
$ws.Range("N4").Value = 78

# Average-price formula for the Miercoles M cell
$ws.Range("G4").Formula = "=(F4*0.4+H4*0.6)/2"

# Mark the week as a "Big Spike" (leading apostrophe forces a text/quote-prefix cell)
$ws.Range("P4").Value = "'Big Spike'"

# Make "Datos" the active sheet with N4 selected, matching the saved view state
$ws.Activate() | Out-Null
$ws.Range("N4").Select() | Out-Null
